$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the marking value for row 11 (B11: 3 -> 5)
$ws.Range("B11").Value = 5

# Update the total correct count for row 12 (B12: 72 -> 120)
$ws.Range("B12").Value = 120

# Update the correct/total marks text (E12: "70/84" -> "120/140")
$ws.Range("E12").Value = "120/140"
